# Excel To CSV 工具修改
# Rename the typed-field headers (string -> String) on both sheets, widen a
# couple of columns, bump the China sheet's sample price to a near-10001
# float, reformat the England sheet's description values into
# comma-separated characters, rename its array-header suffix, and flip
# which sheet/cell is the active selection.

$wb = $excel.ActiveWorkbook
$wsChina   = $wb.Worksheets.Item("China")
$wsEngland = $wb.Worksheets.Item("England")

# --- Header renames (both sheets use the same three renamed headers) ---
$wsChina.Range("B1").Value = "name.String"
$wsChina.Range("C1").Value = "desc.String"
$wsChina.Range("E1").Value = "bornPoint.String.array"

$wsEngland.Range("B1").Value = "name.String"
$wsEngland.Range("C1").Value = "desc.String"
$wsEngland.Range("E1").Value = "bornPoint.float.array.aaa"

# --- China sheet: sample price becomes a near-integer float ---
$wsChina.Range("D2").Value = 10000.999999899999

# --- England sheet: comma-join every description value's characters ---
$wsEngland.Range("C2").Value  = "物,品,描,述1"
$wsEngland.Range("C3").Value  = "物,品,描,述2"
$wsEngland.Range("C4").Value  = "物,品,描,述3"
$wsEngland.Range("C5").Value  = "物,品,描,述4"
$wsEngland.Range("C6").Value  = "物,品,描,述5"
$wsEngland.Range("C7").Value  = "物,品,描,述6"
$wsEngland.Range("C8").Value  = "物,品,描,述7"
$wsEngland.Range("C9").Value  = "物,品,描,述8"
$wsEngland.Range("C10").Value = "物,品,描,述9"

# --- Column width tweaks ---
$wsChina.Columns.Item(1).ColumnWidth = 23.875
$wsEngland.Columns.Item(5).ColumnWidth = 28.25

# --- Selection / active sheet: England's range gets selected first so
#     China (selected last) ends up as the active tab, matching the
#     tabSelected flag moving from England to China. ---
$wsEngland.Select() | Out-Null
$wsEngland.Range("A2:F2").Select() | Out-Null

$wsChina.Select() | Out-Null
$wsChina.Range("A2").Select() | Out-Null

Write-Output "done"
